$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7500298023223877
$ws.Range("B1").Value = 1.407041192054749
$ws.Range("C1").Value = 5.295890808105469
$ws.Range("D1").Value = 3.176616907119751
$ws.Range("E1").Value = 1.522581815719604
